$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Update boolean / numeric input cells
$ws.Range("B28").Value = $false
$ws.Range("B29").Value = $false
$ws.Range("B30").Value = 2009
$ws.Range("B31").Value = $false

# Update formula in C31 to include the extra AND(B31=TRUE, ...) condition
$ws.Range("C31").Formula = '=IF(AND(B31=TRUE,OR(B28<>TRUE,B29<>TRUE)),"demand and profiles must be fix!!!!!!!!!!","- > NOT ACTIVE")'

# Align formatting of B23 and B31 with the rest of the block (same style as B24)
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the view so topLeftCell/selection matches
$ws.Range("C24").Select()

$wb.Save()
